$d = $word.ActiveDocument
$endash = [char]0x2013

# 1. "... предлагаете. <en-dash> Они просто ..." -> hyphen-minus with spaces
$d.Content.Find.Execute(
    ("предлагаете. " + $endash + " Они"),
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "предлагаете. - Они",
    2)

# 2. "... . <en-dash> Почему ты не снабдил их оружием?" -> hyphen-minus
$d.Content.Find.Execute(
    (". " + $endash + " Почему"),
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    ". - Почему",
    2)

# 3. "... замена <en-dash> честно говоря ..." -> hyphen-minus
$d.Content.Find.Execute(
    ("замена " + $endash + " честно"),
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "замена - честно",
    2)

# 4. "-Представьте" -> "- Представьте" (missing space after hyphen)
$d.Content.Find.Execute(
    "-Представьте",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "- Представьте",
    2)
